$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.422.11'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '3.717.18'
$ws.Range('E3').Value = '  +3.36%  '
$origStyle_D4 = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = $origStyle_D4
$ws.Range('E4').Value = '  +0.05%  '
$origStyle_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.46'
$ws.Range('D5').Style = $origStyle_D5
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('E6').Value = '  +7.61%  '
$origStyle_D7 = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '654.79'
$ws.Range('D7').Style = $origStyle_D7
$ws.Range('E7').Value = '  +0.31%  '
$origStyle_D8 = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.422'
$ws.Range('D8').Style = $origStyle_D8
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('E9').Value = '  +1.50%  '
$origStyle_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('D10').Style = $origStyle_D10
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '3.714.54'
$ws.Range('E11').Value = '  +3.37%  '
$origStyle_D12 = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.00'
$ws.Range('D12').Style = $origStyle_D12
$ws.Range('E12').Value = '  +0.42%  '
$origStyle_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.205'
$ws.Range('D13').Style = $origStyle_D13
$ws.Range('E13').Value = '  +0.49%  '
$origStyle_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.85'
$ws.Range('D14').Style = $origStyle_D14
$ws.Range('E14').Value = '  +5.89%  '
$ws.Range('D15').Value = '4.412.44'
$ws.Range('E15').Value = '  +3.48%  '
$origStyle_D16 = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000268'
$ws.Range('D16').Style = $origStyle_D16
$ws.Range('E16').Value = '  +2.78%  '
$ws.Range('D17').Value = '96.368.36'
$ws.Range('E17').Value = '  -0.65%  '
$origStyle_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.99'
$ws.Range('D18').Style = $origStyle_D18
$ws.Range('E18').Value = '  +16.06%  '
$ws.Range('D19').Value = '3.706.25'
$ws.Range('E19').Value = '  +2.81%  '
$origStyle_D20 = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.12'
$ws.Range('D20').Style = $origStyle_D20
$ws.Range('E20').Value = '  +4.44%  '
$origStyle_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.78'
$ws.Range('D21').Style = $origStyle_D21
$ws.Range('E21').Value = '  +0.73%  '
$origStyle_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.526'
$ws.Range('D22').Style = $origStyle_D22
$ws.Range('E22').Value = '  -0.28%  '
$origStyle_D23 = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '523.65'
$ws.Range('D23').Style = $origStyle_D23
$ws.Range('E23').Value = '  +1.10%  '
$origStyle_D24 = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.48'
$ws.Range('D24').Style = $origStyle_D24
$ws.Range('E24').Value = '  +0.07%  '
$origStyle_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.04'
$ws.Range('D25').Style = $origStyle_D25
$origStyle_D26 = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000202'
$ws.Range('D26').Style = $origStyle_D26
$ws.Range('E26').Value = '  -1.97%  '
$origStyle_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.23'
$ws.Range('D27').Style = $origStyle_D27
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('E29').Value = '  -8.01%  '
$origStyle_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.45'
$ws.Range('D30').Style = $origStyle_D30
$ws.Range('E30').Value = '  +3.46%  '
$origStyle_D31 = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.06'
$ws.Range('D31').Style = $origStyle_D31
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  +10.67%  '
$ws.Range('E34').Value = '  -2.20%  '
$origStyle_D35 = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '667.88'
$ws.Range('D35').Style = $origStyle_D35
$ws.Range('E35').Value = '  +8.40%  '
$origStyle_D36 = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '32.78'
$ws.Range('D36').Style = $origStyle_D36
$ws.Range('E36').Value = '  +2.91%  '
$ws.Range('E37').Value = '  +0.29%  '
$origStyle_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.596'
$ws.Range('D38').Style = $origStyle_D38
$ws.Range('E38').Value = '  +2.17%  '
$origStyle_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.87'
$ws.Range('D39').Style = $origStyle_D39
$ws.Range('E39').Value = '  +0.83%  '
$origStyle_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.08'
$ws.Range('D40').Style = $origStyle_D40
$ws.Range('E40').Value = '  +16.22%  '
$origStyle_D41 = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '40.52'
$ws.Range('D41').Style = $origStyle_D41
$ws.Range('E41').Value = '  +23.16%  '
$ws.Range('E42').Value = '  +4.29%  '
$origStyle_D43 = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.978'
$ws.Range('D43').Style = $origStyle_D43
$ws.Range('E43').Value = '  +5.15%  '
$origStyle_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.98'
$ws.Range('D44').Style = $origStyle_D44
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$origStyle_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0460'
$ws.Range('D46').Style = $origStyle_D46
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$origStyle_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.440'
$ws.Range('D47').Style = $origStyle_D47
$ws.Range('E47').Value = '  -3.57%  '
$origStyle_D48 = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.33'
$ws.Range('D48').Style = $origStyle_D48
$ws.Range('E48').Value = '  -0.56%  '
$origStyle_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.61'
$ws.Range('D49').Style = $origStyle_D49
$ws.Range('E49').Value = '  -0.22%  '
$origStyle_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.59'
$ws.Range('D50').Style = $origStyle_D50
$ws.Range('E50').Value = '  -0.93%  '
$ws.Range('E51').Value = '  +2.46%  '
